$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.676.27"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.54%  "
$ws.Range("D3").Value = "'1.742.48"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -5.90%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'235.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -10.43%  "
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").Value = "'0.4926"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -7.68%  "
$ws.Range("D8").Value = "'41.44"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -8.27%  "
$ws.Range("D9").Value = "'0.2551"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -19.36%  "
$ws.Range("D10").Value = "'0.06014"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -13.61%  "
$ws.Range("D11").Value = "'1.748.22"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.34%  "
$ws.Range("D12").Value = "'0.06827"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -12.86%  "
$ws.Range("D13").Value = "'14.77"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -21.80%  "
$ws.Range("D14").Value = "'4.447"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -12.01%  "
$ws.Range("D15").Value = "'76.57"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -14.76%  "
$ws.Range("D16").Value = "'0.5666"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -26.52%  "
$ws.Range("D17").Value = "'1.002"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.07%  "
$ws.Range("D18").Value = "'1.000"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "'25.721.54"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.46%  "
$ws.Range("D20").Value = "'11.26"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -20.39%  "
$ws.Range("D21").Value = "'0.000006566"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -17.61%  "
$ws.Range("D22").Value = "'1.965.71"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.90%  "
$ws.Range("D23").Value = "'4.008"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Value = "'5.049"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -16.33%  "
$ws.Range("D25").Value = "'7.911"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -15.51%  "
$ws.Range("D26").Value = "'137.29"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.16%  "
$ws.Range("D27").Value = "'1.475"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -12.82%  "
$ws.Range("D28").Value = "'1.825"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -17.59%  "
$ws.Range("D29").Value = "'14.66"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -14.40%  "
$ws.Range("D30").Value = "'101.91"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -8.83%  "
$ws.Range("D31").Value = "'3.764"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -12.65%  "
$ws.Range("D32").Value = "'0.07991"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -9.01%  "
$ws.Range("D33").Value = "'3.405"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -17.37%  "
$ws.Range("D34").Value = "'0.04383"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -9.78%  "
$ws.Range("D35").Value = "'0.9998"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("D36").Value = "'2.602"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -9.93%  "
$ws.Range("D37").Value = "'0.9818"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -13.88%  "
$ws.Range("D38").Value = "'0.6016"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -18.55%  "
$ws.Range("D39").Value = "'2.675"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -14.21%  "
$ws.Range("D40").Value = "'1.965"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -16.38%  "
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("D42").Value = "'0.01512"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -12.97%  "
$ws.Range("D43").Value = "'101.80"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.47%  "
$ws.Range("D44").Value = "'0.7567"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -16.68%  "
$ws.Range("D45").Value = "'5.163"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -12.78%  "
$ws.Range("D46").Value = "'0.3748"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -22.50%  "
$ws.Range("D47").Value = "'0.05229"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -10.08%  "
$ws.Range("D48").Value = "'0.1064"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -15.05%  "
$ws.Range("D49").Value = "'30.09"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -14.16%  "
$ws.Range("D50").Value = "'52.20"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -13.74%  "
$ws.Range("D51").Value = "'5.823"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -24.38%  "
